# Update PLC data 2025-10-13 13:44:42
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 157314
$ws.Range("C4").Value = 148375
$ws.Range("C5").Value = 8939
$ws.Range("C8").Value = 63.9
